$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 2 values (text values stay the same, only numeric recompute)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 0.02035233333333333
$ws.Range("N2").Value = 0.061057
$ws.Range("O2").Value = 0.128389957923202
$ws.Range("P2").Value = 0.128389957923202
$ws.Range("S2").Value = 0.128389957923202
$ws.Range("T2").Value = 0.128389957923202

# Update row 3 values
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1381673333333333
$ws.Range("N3").Value = 0.414502
$ws.Range("O3").Value = 0.871610042076798
$ws.Range("P3").Value = 0.871610042076798
$ws.Range("Q3").Value = 0.01617769066955555
$ws.Range("R3").Value = 0.145599216026
$ws.Range("S3").Value = 0.871610042076798
$ws.Range("T3").Value = 0.871610042076798

# Delete rows 4 and 5 (Resolving-Mac rows), shifting rows up
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()
